$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be auto-parsed as numbers
# by Excel, so they remain text cells exactly like the source data (inline strings).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = "60.245.48"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "2.611.23"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "584.87"
$ws.Range("E5").Value = "  +2.59%  "
$ws.Range("D6").Value = "143.13"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("E11").Value = "  +2.05%  "
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").Value = "3.071.45"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").Value = "24.84"
$ws.Range("E14").Value = "  +5.44%  "
$ws.Range("D15").Value = "60.242.92"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("D16").Value = "0.0000140"
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("D17").Value = "2.616.88"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "11.42"
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").Value = "4.63"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").Value = "346.42"
$ws.Range("E20").Value = "  -1.15%  "
$ws.Range("D21").Value = "6.91"
$ws.Range("E21").Value = "  -2.45%  "
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("D23").Value = "0.533"
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("D24").Value = "63.68"
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("D27").Value = "8.02"
$ws.Range("E27").Value = "  +3.64%  "
$ws.Range("E28").Value = "  +5.94%  "
$ws.Range("D29").Value = "0.0₃0797"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "6.49"
$ws.Range("E30").Value = "  +1.96%  "
$ws.Range("D31").Value = "169.16"
$ws.Range("E31").Value = "  +4.45%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("D34").Value = "1.02"
$ws.Range("E34").Value = "  +5.88%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "1.31"
$ws.Range("E35").Value = "  +8.15%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "4.30"
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("D37").Value = "1.63"
$ws.Range("E37").Value = "  +2.32%  "
$ws.Range("D38").Value = "319.12"
$ws.Range("E38").Value = "  +6.44%  "
$ws.Range("D39").Value = "38.36"
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("E40").Value = "  +3.39%  "
$ws.Range("D41").Value = "0.852"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").Value = "135.70"
$ws.Range("E42").Value = "  -2.97%  "
$ws.Range("D43").Value = "0.0991"
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").Value = "19.94"
$ws.Range("E45").Value = "  +1.73%  "
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D47").Value = "5.01"
$ws.Range("E47").Value = "  +2.08%  "
$ws.Range("D48").Value = "0.0550"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "20.00"
$ws.Range("E49").Value = "  +0.92%  "
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("D51").Value = "10.75"
$ws.Range("E51").Value = "  +0.55%  "
